$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.9155336666666667
$ws.Range("H2").Value = 2.746601
$ws.Range("I2").Value = 0.1890240037548773
$ws.Range("J2").Value = 0.1951995261655112
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 111.629865488713
$ws.Range("R2").Value = 1004.668789398417
$ws.Range("S2").Value = 0.04313966587763295
$ws.Range("T2").Value = 0.0472328281646562
$ws.Range("G3").Value = 0.9155336666666667
$ws.Range("H3").Value = 2.746601
$ws.Range("I3").Value = 0.1890240037548773
$ws.Range("J3").Value = 0.1951995261655112
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 135.4205488974433
$ws.Range("R3").Value = 1218.78494007699
$ws.Range("S3").Value = 0.05233364034638249
$ws.Range("T3").Value = 0.05729914201753741
$ws.Range("G4").Value = 0.9155336666666667
$ws.Range("H4").Value = 2.746601
$ws.Range("I4").Value = 0.1890240037548773
$ws.Range("J4").Value = 0.1951995261655112
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 76.45161075863423
$ws.Range("R4").Value = 688.064496827708
$ws.Range("S4").Value = 0.02954493342346456
$ws.Range("T4").Value = 0.03234820518742689
$ws.Range("G5").Value = 0.9155336666666667
$ws.Range("H5").Value = 2.746601
$ws.Range("I5").Value = 0.1890240037548773
$ws.Range("J5").Value = 0.1951995261655112
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 83.37624339032084
$ws.Range("R5").Value = 500.2574603419251
$ws.Range("S5").Value = 0.03222097658403362
$ws.Range("T5").Value = 0.02351877047615459
$ws.Range("G6").Value = 0.9155336666666667
$ws.Range("H6").Value = 2.746601
$ws.Range("I6").Value = 0.1890240037548773
$ws.Range("J6").Value = 0.1951995261655112
$ws.Range("M6").Value = 89.83563
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 82.24754373121
$ws.Range("R6").Value = 740.22789358089
$ws.Range("S6").Value = 0.03178478752336362
$ws.Range("T6").Value = 0.03480058031973606
$ws.Range("I7").Value = 0.7156204889943075
$ws.Range("J7").Value = 0.7390002200311341
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 422.616267460901
$ws.Range("R7").Value = 3803.546407148109
$ws.Range("S7").Value = 0.1633212088261366
$ws.Range("T7").Value = 0.1788173931158901
$ws.Range("I8").Value = 0.7156204889943075
$ws.Range("J8").Value = 0.7390002200311341
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 512.6847252031366
$ws.Range("R8").Value = 4614.16252682823
$ws.Range("S8").Value = 0.1981284098928316
$ws.Range("T8").Value = 0.2169271585354745
$ws.Range("I9").Value = 0.7156204889943075
$ws.Range("J9").Value = 0.7390002200311341
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 289.4359340015018
$ws.Range("R9").Value = 2604.923406013516
$ws.Range("S9").Value = 0.1118533058437475
$ws.Range("T9").Value = 0.1224661310440438
$ws.Range("I10").Value = 0.7156204889943075
$ws.Range("J10").Value = 0.7390002200311341
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 315.6516996797042
$ws.Range("R10").Value = 1893.910198078225
$ws.Range("S10").Value = 0.1219844599675364
$ws.Range("T10").Value = 0.08903903046364482
$ws.Range("I11").Value = 0.7156204889943075
$ws.Range("J11").Value = 0.7390002200311341
$ws.Range("M11").Value = 89.83563
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 311.37858840317
$ws.Range("R11").Value = 2802.40729562853
$ws.Range("S11").Value = 0.1203331044640554
$ws.Range("T11").Value = 0.1317505068720808
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.002153333333333333
$ws.Range("H12").Value = 0.00646
$ws.Range("I12").Value = 0.0004445840747369229
$ws.Range("J12").Value = 0.0004591088909634862
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 0.26255321798
$ws.Range("R12").Value = 2.36297896182
$ws.Range("S12").Value = 0.0001014644069413464
$ws.Range("T12").Value = 0.0001110915163664759
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.002153333333333333
$ws.Range("H13").Value = 0.00646
$ws.Range("I13").Value = 0.0004445840747369229
$ws.Range("J13").Value = 0.0004591088909634862
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 0.3185088572666667
$ws.Range("R13").Value = 2.8665797154
$ws.Range("S13").Value = 0.0001230886163070759
$ws.Range("T13").Value = 0.0001347674662003297
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.002153333333333333
$ws.Range("H14").Value = 0.00646
$ws.Range("I14").Value = 0.0004445840747369229
$ws.Range("J14").Value = 0.0004591088909634862
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 0.1798140339644445
$ws.Range("R14").Value = 1.61832630568
$ws.Range("S14").Value = 0.00006948962368963714
$ws.Range("T14").Value = 0.00007608291321192183
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.002153333333333333
$ws.Range("H15").Value = 0.00646
$ws.Range("I15").Value = 0.0004445840747369229
$ws.Range("J15").Value = 0.0004591088909634862
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 0.1961007559166667
$ws.Range("R15").Value = 1.1766045355
$ws.Range("S15").Value = 0.00007578367179392172
$ws.Range("T15").Value = 0.00005531610061889538
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.002153333333333333
$ws.Range("H16").Value = 0.00646
$ws.Range("I16").Value = 0.0004445840747369229
$ws.Range("J16").Value = 0.0004591088909634862
$ws.Range("M16").Value = 89.83563
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 0.1934460566
$ws.Range("R16").Value = 1.7410145094
$ws.Range("S16").Value = 0.00007475775600494174
$ws.Range("T16").Value = 0.00008185089456586339
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.459699
$ws.Range("H17").Value = 0.9193979999999999
$ws.Range("I17").Value = 0.09491092317607834
$ws.Range("J17").Value = 0.06534114491239122
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 56.050519389561
$ws.Range("R17").Value = 336.303116337366
$ws.Range("S17").Value = 0.0216608760401842
$ws.Range("T17").Value = 0.01581073033503177
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.459699
$ws.Range("H18").Value = 0.9193979999999999
$ws.Range("I18").Value = 0.09491092317607834
$ws.Range("J18").Value = 0.06534114491239122
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 67.99606958666999
$ws.Range("R18").Value = 407.9764175200199
$ws.Range("S18").Value = 0.02627726648347359
$ws.Range("T18").Value = 0.01918033109746915
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.459699
$ws.Range("H19").Value = 0.9193979999999999
$ws.Range("I19").Value = 0.09491092317607834
$ws.Range("J19").Value = 0.06534114491239122
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 38.387150897564
$ws.Range("R19").Value = 230.322905385384
$ws.Range("S19").Value = 0.01483481912716834
$ws.Range("T19").Value = 0.01082824740576076
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.459699
$ws.Range("H20").Value = 0.9193979999999999
$ws.Range("I20").Value = 0.09491092317607834
$ws.Range("J20").Value = 0.06534114491239122
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 41.86408114278751
$ws.Range("R20").Value = 167.45632457115
$ws.Range("S20").Value = 0.01617848830030682
$ws.Range("T20").Value = 0.007872679919011017
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.459699
$ws.Range("H21").Value = 0.9193979999999999
$ws.Range("I21").Value = 0.09491092317607834
$ws.Range("J21").Value = 0.06534114491239122
$ws.Range("M21").Value = 89.83563
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 41.29734927537
$ws.Range("R21").Value = 247.78409565222
$ws.Range("S21").Value = 0.01595947322494538
$ws.Range("T21").Value = 0.01164915615511852
